$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "level": replace the A:D level-definition values.
# Column E (formula meta!$B$1*ROW()) is untouched.
# ---------------------------------------------------------------
$wsLevel = $wb.Worksheets.Item("level")

# Clear out the old A:D values first (rows 1-23 held data in the
# original workbook).
$wsLevel.Range("A1:D23").ClearContents()

$wsLevel.Range("A1").Value = 7

$wsLevel.Range("B2").Value = 8
$wsLevel.Range("D2").Value = 10

$wsLevel.Range("C3").Value = 9

$wsLevel.Range("B5").Value = 11

$wsLevel.Range("A6").Value = 9
$wsLevel.Range("D6").Value = 10

$wsLevel.Range("C7").Value = 8

$wsLevel.Range("B9").Value = 7

$wsLevel.Range("A11").Value = 10

$wsLevel.Range("B13").Value = 7
$wsLevel.Range("D13").Value = 11

# ---------------------------------------------------------------
# Sheet "enemies": update existing rows 2-7 and append new rows 8-12.
# ---------------------------------------------------------------
$wsEnemies = $wb.Worksheets.Item("enemies")

$wsEnemies.Range("A2").Value = 1
$wsEnemies.Range("B2").Value = 0
$wsEnemies.Range("C2").Value = 0
$wsEnemies.Range("D2").Value = 2
$wsEnemies.Range("E2").Value = 1
$wsEnemies.Range("F2").Value = 'EnemyPrefabs/Special Enemies/Halloween Bee/Halloween Bee'

$wsEnemies.Range("A3").Value = 2
$wsEnemies.Range("B3").Value = 0
$wsEnemies.Range("C3").Value = 0
$wsEnemies.Range("D3").Value = 2
$wsEnemies.Range("E3").Value = 1
$wsEnemies.Range("F3").Value = 'EnemyPrefabs/Special Enemies/Steampunk Fly/Steampunk Fly'

$wsEnemies.Range("A4").Value = 3
$wsEnemies.Range("B4").Value = 0
$wsEnemies.Range("C4").Value = 0
$wsEnemies.Range("D4").Value = 2
$wsEnemies.Range("E4").Value = 1
$wsEnemies.Range("F4").Value = 'EnemyPrefabs/Special Enemies/Bionic Lady Bird/Bionic Lady Bird'

$wsEnemies.Range("A5").Value = 4
$wsEnemies.Range("B5").Value = 0
$wsEnemies.Range("C5").Value = 0
$wsEnemies.Range("D5").Value = 3
$wsEnemies.Range("E5").Value = 1.5
$wsEnemies.Range("F5").Value = 'EnemyPrefabs/Special Enemies/Halloween Bee/Halloween Bee'

$wsEnemies.Range("A6").Value = 5
$wsEnemies.Range("B6").Value = 0
$wsEnemies.Range("C6").Value = 0
$wsEnemies.Range("D6").Value = 3
$wsEnemies.Range("E6").Value = 1.5
$wsEnemies.Range("F6").Value = 'EnemyPrefabs/Special Enemies/Steampunk Fly/Steampunk Fly'
$wsEnemies.Range("G6").Value = 'EnemyPrefabs/Bullet Enemies//'

$wsEnemies.Range("A7").Value = 6
$wsEnemies.Range("B7").Value = 0
$wsEnemies.Range("C7").Value = 0
$wsEnemies.Range("D7").Value = 3
$wsEnemies.Range("E7").Value = 1.5
$wsEnemies.Range("F7").Value = 'EnemyPrefabs/Special Enemies/Bionic Lady Bird/Bionic Lady Bird'
$wsEnemies.Range("G7").Value = 'EnemyPrefabs/Special Enemies//'

$wsEnemies.Range("A8").Value = 7
$wsEnemies.Range("B8").Value = 4
$wsEnemies.Range("C8").Value = 0
$wsEnemies.Range("D8").Value = 0
$wsEnemies.Range("E8").Value = 1
$wsEnemies.Range("F8").Value = 'EnemyPrefabs/Arrow Enemies/Bee/Bee Arrow'
$wsEnemies.Range("A8").HorizontalAlignment = -4108
$wsEnemies.Range("A8").VerticalAlignment = -4108

$wsEnemies.Range("A9").Value = 8
$wsEnemies.Range("B9").Value = 0
$wsEnemies.Range("C9").Value = 2
$wsEnemies.Range("D9").Value = 0
$wsEnemies.Range("E9").Value = 1
$wsEnemies.Range("F9").Value = 'EnemyPrefabs/Bullet Enemies/Neo Fly/Neo Fly'
$wsEnemies.Range("A9").HorizontalAlignment = -4108
$wsEnemies.Range("A9").VerticalAlignment = -4108

$wsEnemies.Range("A10").Value = 9
$wsEnemies.Range("B10").Value = 0
$wsEnemies.Range("C10").Value = 0
$wsEnemies.Range("D10").Value = 2
$wsEnemies.Range("E10").Value = 1
$wsEnemies.Range("F10").Value = 'EnemyPrefabs/Special Enemies/Steampunk Fly/Steampunk Fly'
$wsEnemies.Range("A10").HorizontalAlignment = -4108
$wsEnemies.Range("A10").VerticalAlignment = -4108

$wsEnemies.Range("A11").Value = 10
$wsEnemies.Range("B11").Value = 0
$wsEnemies.Range("C11").Value = 0
$wsEnemies.Range("D11").Value = 1
$wsEnemies.Range("E11").Value = 2.5
$wsEnemies.Range("F11").Value = 'EnemyPrefabs/Special Enemies/Halloween Bee/Halloween Bee'
$wsEnemies.Range("A11").HorizontalAlignment = -4108
$wsEnemies.Range("A11").VerticalAlignment = -4108

$wsEnemies.Range("A12").Value = 11
$wsEnemies.Range("B12").Value = 2
$wsEnemies.Range("C12").Value = 0
$wsEnemies.Range("D12").Value = 0
$wsEnemies.Range("E12").Value = 1
$wsEnemies.Range("F12").Value = 'EnemyPrefabs/Arrow Enemies/Fly/Fly Arrow'
$wsEnemies.Range("A12").HorizontalAlignment = -4108
$wsEnemies.Range("A12").VerticalAlignment = -4108

# ---------------------------------------------------------------
# Sheet "misc": refresh the full reference list of enemy prefab
# paths (now includes the new Butterfly Arrow entry + 2 placeholders).
# ---------------------------------------------------------------
$wsMisc = $wb.Worksheets.Item("misc")

$wsMisc.Range("A1").Value = 'all enemies:'
$wsMisc.Range("A2").Value = 'EnemyPrefabs/Arrow Enemies/Bee/Bee Arrow'
$wsMisc.Range("A3").Value = 'EnemyPrefabs/Bullet Enemies/Funky Bee/Funky Bee'
$wsMisc.Range("A4").Value = 'EnemyPrefabs/Special Enemies/Halloween Bee/Halloween Bee'
$wsMisc.Range("A5").Value = 'EnemyPrefabs/Arrow Enemies/Butterfly/Butterfly Arrow'
$wsMisc.Range("A6").Value = 'EnemyPrefabs/Special Enemies/Halloween Bee/Halloween Bee'
$wsMisc.Range("A7").Value = 'EnemyPrefabs/Special Enemies/Halloween Bee/Halloween Bee'
$wsMisc.Range("A8").Value = 'EnemyPrefabs/Arrow Enemies/Dragonfly/Dragonfly Arrow'
$wsMisc.Range("A9").Value = 'EnemyPrefabs/Special Enemies/Halloween Bee/Halloween Bee'
$wsMisc.Range("A10").Value = 'EnemyPrefabs/Special Enemies/Halloween Bee/Halloween Bee'
$wsMisc.Range("A11").Value = 'EnemyPrefabs/Arrow Enemies/Fly/Fly Arrow'
$wsMisc.Range("A12").Value = 'EnemyPrefabs/Bullet Enemies/Neo Fly/Neo Fly'
$wsMisc.Range("A13").Value = 'EnemyPrefabs/Special Enemies/Steampunk Fly/Steampunk Fly'
$wsMisc.Range("A14").Value = 'EnemyPrefabs/Arrow Enemies/Lady Bird/Lady Bird Arrow'
$wsMisc.Range("A15").Value = 'EnemyPrefabs/Bullet Enemies/Magic Lady Bird/Magic Lady Bird'
$wsMisc.Range("A16").Value = 'EnemyPrefabs/Special Enemies/Bionic Lady Bird/Bionic Lady Bird'
